$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing style applied to J3 (a custom font colour) needs to move to
# the new "binary[binary]" header cell (L1). Upgrade the font on J3 first -
# this updates the existing style slot in place - then copy that formatted
# cell over to L1 before overwriting the values, so the destination cell
# picks up the same style slot and J3 can be cleared back to the default.
$ws.Range("J3").Font.Color = 0
$ws.Range("J3").Copy($ws.Range("L1"))
$ws.Range("J3").ClearFormats()
$ws.Range("J3").Value = 14.3

# New column headers for the Edm.Byte / Edm.Binary special-case columns
$ws.Range("K1").Value = "byte[byte]"
$ws.Range("L1").Value = "binary[binary]"

# New column data for the two sample rows
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = "dGhpc2lzYXRlc3Q="
$ws.Range("K3").Value = 4
$ws.Range("L3").Value = "dGhpc2lzYXRlc3Q="

# Update the active selection to the new last cell
$ws.Range("L3").Select()
